$wb = $excel.ActiveWorkbook

$updates = @(
    @{Sheet="展览"; Cell="F2"; Value=266},
    @{Sheet="展览"; Cell="F3"; Value=564},
    @{Sheet="展览"; Cell="F6"; Value=1073},
    @{Sheet="展览"; Cell="F7"; Value=1404},
    @{Sheet="展览"; Cell="F9"; Value=100},
    @{Sheet="展览"; Cell="F10"; Value=736},
    @{Sheet="展览"; Cell="F13"; Value=118},
    @{Sheet="展览"; Cell="F14"; Value=411},
    @{Sheet="展览"; Cell="F15"; Value=1303},
    @{Sheet="展览"; Cell="F16"; Value=98},
    @{Sheet="展览"; Cell="F18"; Value=264},
    @{Sheet="展览"; Cell="F19"; Value=5215},
    @{Sheet="展览"; Cell="F22"; Value=194},
    @{Sheet="展览"; Cell="F23"; Value=9},
    @{Sheet="展览"; Cell="F24"; Value=5593},
    @{Sheet="展览"; Cell="F25"; Value=52},
    @{Sheet="展览"; Cell="F26"; Value=115},
    @{Sheet="展览"; Cell="F27"; Value=87},
    @{Sheet="展览"; Cell="F29"; Value=14095},
    @{Sheet="展览"; Cell="F30"; Value=1408},
    @{Sheet="展览"; Cell="F31"; Value=191},
    @{Sheet="展览"; Cell="F32"; Value=88},
    @{Sheet="展览"; Cell="F34"; Value=415},
    @{Sheet="展览"; Cell="F35"; Value=571},
    @{Sheet="展览"; Cell="F36"; Value=4159},
    @{Sheet="展览"; Cell="F37"; Value=104},

    @{Sheet="全部类型"; Cell="F2"; Value=266},
    @{Sheet="全部类型"; Cell="F3"; Value=564},
    @{Sheet="全部类型"; Cell="F6"; Value=1073},
    @{Sheet="全部类型"; Cell="F7"; Value=1404},
    @{Sheet="全部类型"; Cell="F9"; Value=100},
    @{Sheet="全部类型"; Cell="F10"; Value=736},
    @{Sheet="全部类型"; Cell="F13"; Value=118},
    @{Sheet="全部类型"; Cell="F14"; Value=411},
    @{Sheet="全部类型"; Cell="F15"; Value=1303},
    @{Sheet="全部类型"; Cell="F16"; Value=98},
    @{Sheet="全部类型"; Cell="F18"; Value=264},
    @{Sheet="全部类型"; Cell="F20"; Value=5215},
    @{Sheet="全部类型"; Cell="F24"; Value=194},
    @{Sheet="全部类型"; Cell="F25"; Value=9},
    @{Sheet="全部类型"; Cell="F27"; Value=5593},
    @{Sheet="全部类型"; Cell="F28"; Value=52},
    @{Sheet="全部类型"; Cell="F29"; Value=115},
    @{Sheet="全部类型"; Cell="F30"; Value=87},
    @{Sheet="全部类型"; Cell="F32"; Value=14095},
    @{Sheet="全部类型"; Cell="F33"; Value=1408},
    @{Sheet="全部类型"; Cell="F34"; Value=191},
    @{Sheet="全部类型"; Cell="F35"; Value=88},
    @{Sheet="全部类型"; Cell="F37"; Value=415},
    @{Sheet="全部类型"; Cell="F38"; Value=571},
    @{Sheet="全部类型"; Cell="F39"; Value=4159},
    @{Sheet="全部类型"; Cell="F40"; Value=104}
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
